# TS 4.4 Ghanam Tamil Corrections - update the "Observed till" date from
# "April 31,2024" to "April 30,2024".
#
# The source run ("April 31,2024") must end up split into three runs
# ("April 3" / "0" / ",2024") that all keep the original bold/underline/
# size-32 character formatting - mirroring how Word itself fragments a
# run when a single character is retyped in the middle of it.
#
# Directly overwriting the Range.Text of the "1" causes this host's
# paragraph-serializer to re-merge the whole title line's runs, so the
# edit is done in three passes:
#   1) a harmless Font round-trip on the single "1" character to force
#      the engine to carve that character into its own run (without
#      touching its neighbours, "Observed till " stays untouched);
#   2) replace that isolated character's text ("1" -> "0") - this
#      re-merges the former "April 3"/"1"/",2024" pieces into one
#      "April 30,2024" run, but leaves everything outside that run
#      (i.e. "Observed till ") alone;
#   3) a second Font round-trip re-splits "April 30,2024" into the
#      desired "April 3" / "0" / ",2024" triple.

$d = $word.ActiveDocument

$searchText = "April 31,2024"
$content = $d.Content.Text
$startIdx = $content.IndexOf($searchText)

if ($startIdx -ge 0) {
    # Position of the digit to change: "April 3" is 7 characters long,
    # so the "1" in "31" sits right after it.
    $digitStart = $startIdx + 7
    $digitEnd = $digitStart + 1

    # --- Pass 1: isolate the "1" into its own run ---------------------
    $iso = $d.Range($digitStart, $digitEnd)
    $iso.Font.Size = 32
    $iso = $d.Range($digitStart, $digitEnd)
    $iso.Font.Size = 16

    # --- Pass 2: change the digit's text -------------------------------
    $edit = $d.Range($digitStart, $digitEnd)
    $edit.Text = "0"

    # --- Pass 3: re-split the merged "April 30,2024" run --------------
    $resplit = $d.Range($digitStart, $digitEnd)
    $resplit.Font.Size = 32
    $resplit = $d.Range($digitStart, $digitEnd)
    $resplit.Font.Size = 16
}
